$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10089.2
$ws.Range("I6").Value = 10089.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 30267.6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -30155.6
$ws.Range("N6").ClearContents()

$ws.Range("H12").Value = 179.6
$ws.Range("I12").Value = 166
$ws.Range("K12").Value = 166
$ws.Range("M12").Value = 4

$ws.Range("H21").Value = 7964
$ws.Range("I21").Value = 8902.125
$ws.Range("J21").Value = 6463
$ws.Range("K21").Value = 8902.125
$ws.Range("L21").Value = 6463
$ws.Range("M21").Value = -8434.125
$ws.Range("N21").Value = -7399

$ws.Range("H23").Value = 7964
$ws.Range("I23").Value = 8902.125
$ws.Range("J23").Value = 6463
$ws.Range("K23").Value = 8902.125
$ws.Range("L23").Value = 6463
$ws.Range("M23").Value = -8668.125
$ws.Range("N23").Value = -6931

$ws.Range("H29").Value = 1877.2222
$ws.Range("J29").Value = 2732.5
$ws.Range("L29").Value = 8197.5
$ws.Range("N29").Value = -8759.5

$ws.Range("H38").Value = 1277.3684
$ws.Range("I38").Value = 88.3
$ws.Range("J38").Value = 2598.5557
$ws.Range("K38").Value = 264.9
$ws.Range("L38").Value = 7795.6671
$ws.Range("M38").Value = 107.1
$ws.Range("N38").Value = -8539.667099999999

$ws.Range("H58").Value = 723.8261
$ws.Range("J58").Value = 1024.1818
$ws.Range("L58").Value = 3072.5454
$ws.Range("N58").Value = -3372.5454

$ws.Range("H62").Value = 12347206
$ws.Range("I62").Value = 22222982
$ws.Range("J62").Value = 2486.5
$ws.Range("K62").Value = 22222982
$ws.Range("L62").Value = 2486.5
$ws.Range("M62").Value = -22222358
$ws.Range("N62").Value = -3734.5

$ws.Range("H65").Value = 12347206
$ws.Range("I65").Value = 22222982
$ws.Range("J65").Value = 2486.5
$ws.Range("K65").Value = 111114910
$ws.Range("L65").Value = 12432.5
$ws.Range("M65").Value = -111111790
$ws.Range("N65").Value = -18672.5

$ws.Range("H76").Value = 5363.4546
$ws.Range("I76").Value = 5800
$ws.Range("K76").Value = 5800
$ws.Range("M76").Value = -5485

$ws.Range("H79").Value = 5363.4546
$ws.Range("I79").Value = 5800
$ws.Range("K79").Value = 5800
$ws.Range("M79").Value = -4708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7007.242
$ws.Range("I32").Value = 5816.5127
$ws.Range("J32").Value = 12470.588
$ws.Range("K32").Value = 5816.5127
$ws.Range("L32").Value = 12470.588
$ws.Range("M32").Value = -5529.5127
$ws.Range("N32").Value = -13044.588

$ws.Range("H74").Value = 2202.3914
$ws.Range("I74").Value = 1461.0834
$ws.Range("J74").Value = 3011.0908
$ws.Range("K74").Value = 1461.0834
$ws.Range("L74").Value = 3011.0908
$ws.Range("M74").Value = -587.0834
$ws.Range("N74").Value = -4759.0908

$ws.Range("H77").Value = 2202.3914
$ws.Range("I77").Value = 1461.0834
$ws.Range("J77").Value = 3011.0908
$ws.Range("K77").Value = 7305.416999999999
$ws.Range("L77").Value = 15055.454
$ws.Range("M77").Value = -2937.416999999999
$ws.Range("N77").Value = -23791.454

$ws.Range("H132").Value = 2989.2285
$ws.Range("I132").Value = 2792.8823
$ws.Range("J132").Value = 3174.6667
$ws.Range("K132").Value = 8378.6469
$ws.Range("L132").Value = 9524.000100000001
$ws.Range("M132").Value = -5848.6469
$ws.Range("N132").Value = -14584.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 20542
$ws.Range("J27").Value = 20542
$ws.Range("L27").Value = 20542
$ws.Range("N27").Value = -20926

$ws.Range("H29").Value = 1457.5
$ws.Range("I29").Value = 1457.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1457.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1168.5
$ws.Range("N29").ClearContents()

$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1466

$ws.Range("H75").Value = 10720.77
$ws.Range("I75").Value = 4355.6665
$ws.Range("J75").Value = 16176.571
$ws.Range("K75").Value = 4355.6665
$ws.Range("L75").Value = 16176.571
$ws.Range("M75").Value = -3419.6665
$ws.Range("N75").Value = -18048.571

$ws.Range("H78").Value = 10720.77
$ws.Range("I78").Value = 4355.6665
$ws.Range("J78").Value = 16176.571
$ws.Range("K78").Value = 13066.9995
$ws.Range("L78").Value = 48529.713
$ws.Range("M78").Value = -8386.999500000002
$ws.Range("N78").Value = -57889.713

$ws.Range("H94").Value = 8928885
$ws.Range("I94").Value = 9615706
$ws.Range("J94").Value = 209
$ws.Range("K94").Value = 9615706
$ws.Range("L94").Value = 209
$ws.Range("M94").Value = -9615255
$ws.Range("N94").Value = -1111

$ws.Range("H134").Value = 1406.3334
$ws.Range("I134").Value = 1261.4546
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3784.3638
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1249.3638
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4574.9165
$ws.Range("I58").Value = 1273.4445
$ws.Range("J58").Value = 7876.3887
$ws.Range("K58").Value = 1273.4445
$ws.Range("L58").Value = 7876.3887
$ws.Range("M58").Value = -1070.4445
$ws.Range("N58").Value = -8282.3887

$ws.Range("H95").Value = 24731.572
$ws.Range("J95").Value = 24731.572
$ws.Range("L95").Value = 24731.572
$ws.Range("N95").Value = -30223.572

$ws.Range("H99").Value = 1576.3334
$ws.Range("I99").Value = 1476.8182
$ws.Range("J99").Value = 1850
$ws.Range("K99").Value = 1476.8182
$ws.Range("L99").Value = 1850
$ws.Range("M99").Value = 21.18180000000007
$ws.Range("N99").Value = -4846

$ws.Range("H126").Value = 1576.3334
$ws.Range("I126").Value = 1476.8182
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 4430.4546
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -1960.4546
$ws.Range("N126").Value = -10490

$ws.Range("H134").Value = 16130717
$ws.Range("I134").Value = 1525.6818
$ws.Range("J134").Value = 55557628
$ws.Range("K134").Value = 4577.0454
$ws.Range("L134").Value = 166672884
$ws.Range("M134").Value = -2042.0454
$ws.Range("N134").Value = -166677954

$ws.Range("H136").Value = 4574.9165
$ws.Range("I136").Value = 1273.4445
$ws.Range("J136").Value = 7876.3887
$ws.Range("K136").Value = 3820.3335
$ws.Range("L136").Value = 23629.1661
$ws.Range("M136").Value = -1270.3335
$ws.Range("N136").Value = -28729.1661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2649.5
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H80").Value = 3681.0557
$ws.Range("I80").Value = 6000
$ws.Range("J80").Value = 3544.647
$ws.Range("K80").Value = 18000
$ws.Range("L80").Value = 10633.941
$ws.Range("M80").Value = -17064
$ws.Range("N80").Value = -12505.941

$ws.Range("H83").Value = 3681.0557
$ws.Range("I83").Value = 6000
$ws.Range("J83").Value = 3544.647
$ws.Range("K83").Value = 54000
$ws.Range("L83").Value = 31901.823
$ws.Range("M83").Value = -49320
$ws.Range("N83").Value = -41261.823

$ws.Range("H95").Value = 22083.166
$ws.Range("J95").Value = 22083.166
$ws.Range("L95").Value = 66249.49800000001
$ws.Range("N95").Value = -70367.49800000001

$ws.Range("H122").Value = 2333.7778
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 2757.7144
$ws.Range("K122").Value = 7650
$ws.Range("L122").Value = 24819.4296
$ws.Range("M122").Value = -5200
$ws.Range("N122").Value = -29719.4296

$ws.Range("H131").Value = 26319234
$ws.Range("J131").Value = 4650.963
$ws.Range("L131").Value = 13952.889
$ws.Range("N131").Value = -24032.889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3645.4546
$ws.Range("I80").Value = 1840
$ws.Range("K80").Value = 1840
$ws.Range("M80").Value = -842

$ws.Range("H83").Value = 3645.4546
$ws.Range("I83").Value = 1840
$ws.Range("K83").Value = 9200
$ws.Range("M83").Value = -4208

$ws.Range("H92").Value = 29440.2
$ws.Range("J92").Value = 29440.2
$ws.Range("L92").Value = 29440.2
$ws.Range("N92").Value = -33184.2

$ws.Range("H126").Value = 1793.6364
$ws.Range("J126").Value = 2450.2856
$ws.Range("L126").Value = 7350.8568
$ws.Range("N126").Value = -12290.8568

$ws.Range("H132").Value = 3016.5
$ws.Range("I132").Value = 3134.2942
$ws.Range("K132").Value = 9402.882599999999
$ws.Range("M132").Value = -6872.882599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 3276.6667
$ws.Range("J94").Value = 3276.6667
$ws.Range("L94").Value = 3276.6667
$ws.Range("N94").Value = -4628.6667

$ws.Range("H130").Value = 41107.25
$ws.Range("J130").Value = 41107.25
$ws.Range("L130").Value = 41107.25
$ws.Range("N130").Value = -51147.25

$ws.Range("H136").Value = 2520
$ws.Range("I136").Value = 1900
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 5700
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -3150
$ws.Range("N136").Value = -13899.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2994.8948
$ws.Range("I132").Value = 3108.3333
$ws.Range("J132").Value = 3108.3333
$ws.Range("K132").Value = 9324.999899999999
$ws.Range("M132").Value = -6794.999899999999
